$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6846
$ws1.Range("F4").Value = 113
$ws1.Range("F7").Value = 82
$ws1.Range("F8").Value = 591

# Sheet "全部类型" (fourth sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6846
$ws4.Range("F5").Value = 113
$ws4.Range("F9").Value = 82
$ws4.Range("F10").Value = 591
